$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New calibration data for rows 1-108 (columns A:D)
$rows = @(
    @(872.73,1030.79,2499.4,0),
    @(873.53499999999997,1804.723,2499.4,-16.079999999999998),
    @(871.07899999999995,1698.8040000000001,2499.4,-13.98),
    @(869.64599999999996,1599.722,2499.4,-11.95),
    @(868.97299999999996,1503.66,2499.4,-9.9499999999999993),
    @(868.48099999999999,1408.8130000000001,2499.4,-7.98),
    @(868.66399999999999,1313.432,2499.4,-5.95),
    @(869.46299999999997,1219.778,2499.4,-3.98),
    @(870.92200000000003,1125.5039999999999,2499.4,-1.98),
    @(872.68600000000004,1032.056,2499.4,0),
    @(875.46,936.33199999999999,2499.4,2.02),
    @(878.47500000000002,842.59100000000001,2499.4,3.97),
    @(882.27,745.21,2499.4,6),
    @(886.57799999999997,648.93799999999999,2499.4,8),
    @(891.74,552.02700000000004,2499.4,9.9700000000000006),
    @(897.26900000000001,452.77199999999999,2499.4,12),
    @(903.63099999999997,353.27300000000002,2499.4,13.97),
    @(910.80700000000002,250.10499999999999,2499.4,16),
    @(853.92399999999998,1030.81,2599.1,0),
    @(853.88900000000001,1803.374,2599.1,-16.079999999999998),
    @(851.88499999999999,1698.5909999999999,2599.1,-13.98),
    @(850.39800000000002,1599.5519999999999,2599.1,-11.95),
    @(849.78300000000002,1503.2909999999999,2599.1,-9.9499999999999993),
    @(849.32799999999997,1408.096,2599.1,-7.95),
    @(849.73500000000001,1314.393,2599.1,-5.98),
    @(850.54499999999996,1220.135,2599.1,-3.98),
    @(851.96299999999997,1126.201,2599.1,-1.98),
    @(853.93200000000002,1031.4390000000001,2599.1,0),
    @(856.45299999999997,936.33699999999999,2599.1,2.02),
    @(859.48400000000004,842.64300000000003,2599.1,4),
    @(863.19399999999996,747.01599999999996,2599.1,5.97),
    @(867.51700000000005,649.33699999999999,2599.1,8),
    @(872.505,552.41399999999999,2599.1,10),
    @(877.971,454.13,2599.1,12),
    @(884.18200000000002,354.6,2599.1,13.95),
    @(891.06500000000005,251.108,2599.1,16),
    @(836.202,1030.6759999999999,2698.8,0),
    @(835.70399999999995,1800.828,2698.8,-16.079999999999998),
    @(833.89400000000001,1696.4649999999999,2698.8,-13.98),
    @(832.63300000000004,1597.44,2698.8,-11.95),
    @(831.99900000000002,1501.7270000000001,2698.8,-9.9499999999999993),
    @(831.89800000000002,1406.905,2698.8,-7.95),
    @(832.03700000000003,1313.171,2698.8,-5.98),
    @(833,1218.6769999999999,2698.8,-3.98),
    @(834.71600000000001,1124.93,2698.8,-1.98),
    @(836.34799999999996,1029.9559999999999,2698.8,0),
    @(838.98099999999999,935.47900000000004,2698.8,2.02),
    @(841.947,841.91899999999998,2698.8,4),
    @(845.55,745.94500000000005,2698.8,6),
    @(849.91,649.75099999999998,2698.8,8),
    @(854.82799999999997,553.55499999999995,2698.8,9.98),
    @(860.17200000000003,453.38900000000001,2698.8,12),
    @(866.40899999999999,354.62900000000002,2698.8,13.98),
    @(873.00800000000004,253.001,2698.8,16),
    @(820.21600000000001,1030.404,2798.5,0),
    @(819.03300000000002,1800.498,2798.5,-16.079999999999998),
    @(817.04100000000005,1695.5050000000001,2798.5,-13.98),
    @(816.00599999999997,1597.6690000000001,2798.5,-11.95),
    @(815.81399999999996,1501.34,2798.5,-9.9499999999999993),
    @(815.72500000000002,1406.44,2798.5,-7.95),
    @(815.995,1312.652,2798.5,-5.98),
    @(816.99599999999998,1218.4960000000001,2798.5,-3.98),
    @(818.572,1125.4780000000001,2798.5,-1.98),
    @(820.14200000000005,1030.741,2798.5,0),
    @(822.91899999999998,935.82899999999995,2798.5,2.02),
    @(825.83299999999997,842.053,2798.5,4),
    @(829.50699999999995,746.23199999999997,2798.5,6),
    @(833.65300000000002,650.75400000000002,2798.5,7.98),
    @(838.46600000000001,553.33399999999995,2798.5,10),
    @(843.81899999999996,455.34,2798.5,12),
    @(849.58900000000006,355.87200000000001,2798.5,13.98),
    @(856.22,253.75,2798.5,16),
    @(805.11800000000005,1030.8679999999999,2898.2,0),
    @(803.49900000000002,1799.8150000000001,2898.2,-16.079999999999998),
    @(801.95399999999995,1695.347,2898.2,-13.98),
    @(800.976,1597.4839999999999,2898.2,-11.95),
    @(800.279,1501.5119999999999,2898.2,-9.9499999999999993),
    @(800.21900000000005,1406.653,2898.2,-7.95),
    @(800.98,1312.6130000000001,2898.2,-5.98),
    @(801.98900000000003,1218.5150000000001,2898.2,-3.98),
    @(803.51800000000003,1124.5730000000001,2898.2,-1.98),
    @(805.11900000000003,1030.896,2898.2,0),
    @(807.88,936.10400000000004,2898.2,2.02),
    @(810.87199999999996,842.62800000000004,2898.2,4),
    @(814.40800000000002,747.76400000000001,2898.2,6),
    @(818.46400000000006,651.37400000000002,2898.2,7.98),
    @(823.14599999999996,553.80200000000002,2898.2,10),
    @(828.31399999999996,455.68400000000003,2898.2,12),
    @(834.29300000000001,356.904,2898.2,13.98),
    @(840.53,254.417,2898.2,16),
    @(791.07299999999998,1030.72,2997.8,0),
    @(788.56200000000001,1798.2139999999999,2997.8,-16.079999999999998),
    @(787.10400000000004,1692.9,2997.8,-13.98),
    @(786.245,1595.3209999999999,2997.8,-11.95),
    @(786.03700000000003,1499.5889999999999,2997.8,-9.9499999999999993),
    @(786.03399999999999,1405.4739999999999,2997.8,-7.98),
    @(786.86599999999999,1310.2539999999999,2997.8,-5.95),
    @(787.88300000000004,1217.5329999999999,2997.8,-3.98),
    @(789.26900000000001,1123.529,2997.8,-1.98),
    @(791.06399999999996,1030.019,2997.8,0),
    @(793.85799999999995,935.255,2997.8,2.02),
    @(796.73400000000004,841.83600000000001,2997.8,4),
    @(800.29899999999998,746.55600000000004,2997.8,6),
    @(804.452,650.74400000000003,2997.8,7.98),
    @(809.04300000000001,552.47,2997.8,10),
    @(814.08500000000004,455.19299999999998,2997.8,12),
    @(819.64400000000001,356.23,2997.8,13.98),
    @(825.98400000000004,254.01900000000001,2997.8,16)
)

$n = $rows.Count
$arr = New-Object 'object[,]' $n,4
for ($i = 0; $i -lt $n; $i++) {
    for ($j = 0; $j -lt 4; $j++) {
        $arr[$i,$j] = $rows[$i][$j]
    }
}
$ws.Range("A1:D108").Value = $arr

# Remove the now-unused trailing rows (250-267) so the sheet's used range
# shrinks from A1:D267 down to A1:D249
$ws.Range("A250:D267").EntireRow.Delete() | Out-Null

# Update the visible selection to match the author's final view
$ws.Range("A1:D108").Select() | Out-Null
